$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 11252.75
$ws.Range("I21").Value = 10337.333
$ws.Range("J21").Value = 13999
$ws.Range("K21").Value = 10337.333
$ws.Range("L21").Value = 13999
$ws.Range("M21").Value = -9869.333000000001
$ws.Range("N21").Value = -14935

$ws.Range("H23").Value = 11252.75
$ws.Range("I23").Value = 10337.333
$ws.Range("J23").Value = 13999
$ws.Range("K23").Value = 10337.333
$ws.Range("L23").Value = 13999
$ws.Range("M23").Value = -10103.333
$ws.Range("N23").Value = -14467

$ws.Range("H33").Value = 5727.3076
$ws.Range("J33").Value = 3198.2
$ws.Range("L33").Value = 3198.2
$ws.Range("N33").Value = -3656.2

$ws.Range("H39").Value = 1185.5
$ws.Range("I39").Value = 443.0909
$ws.Range("J39").Value = 2818.8
$ws.Range("K39").Value = 1329.2727
$ws.Range("L39").Value = 8456.400000000001
$ws.Range("M39").Value = -1033.2727
$ws.Range("N39").Value = -9048.400000000001

$ws.Range("H62").Value = 1445.5714
$ws.Range("I62").Value = 1445.5714
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1445.5714
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -821.5714
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 1445.5714
$ws.Range("I65").Value = 1445.5714
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 7227.857
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -4107.857
$ws.Range("N65").ClearContents()

$ws.Range("H132").Value = 2813.8223
$ws.Range("I132").Value = 2836.8635
$ws.Range("K132").Value = 8510.5905
$ws.Range("M132").Value = -5980.5905

$ws.Range("H141").Value = 516.46155
$ws.Range("I141").Value = 516.46155
$ws.Range("K141").Value = 1549.38465
$ws.Range("M141").Value = 3630.61535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 7298.6665
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 7298.6665
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 7298.6665
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -7924.6665

$ws.Range("H61").Value = 4989.1816
$ws.Range("I61").Value = 3547.9023
$ws.Range("K61").Value = 3547.9023
$ws.Range("M61").Value = -3335.9023

$ws.Range("H97").Value = 606.3333
$ws.Range("I97").Value = 606.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 606.3333
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -110.3333
$ws.Range("N97").ClearContents()

$ws.Range("H136").Value = 4989.1816
$ws.Range("I136").Value = 3547.9023
$ws.Range("K136").Value = 10643.7069
$ws.Range("M136").Value = -8093.706900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 15715.857
$ws.Range("I3").Value = 12004.333
$ws.Range("K3").Value = 12004.333
$ws.Range("M3").Value = -11891.333

$ws.Range("H31").Value = 5863.844
$ws.Range("J31").Value = 4919.5
$ws.Range("L31").Value = 4919.5
$ws.Range("N31").Value = -5509.5

$ws.Range("H34").Value = 5863.844
$ws.Range("J34").Value = 4919.5
$ws.Range("L34").Value = 4919.5
$ws.Range("N34").Value = -5323.5

$ws.Range("H86").Value = 19612906
$ws.Range("J86").Value = 7587
$ws.Range("L86").Value = 7587
$ws.Range("N86").Value = -9833

$ws.Range("H89").Value = 19612906
$ws.Range("J89").Value = 7587
$ws.Range("L89").Value = 37935
$ws.Range("N89").Value = -49167

$ws.Range("H132").Value = 1621.0834
$ws.Range("I132").Value = 1601.7
$ws.Range("K132").Value = 4805.1
$ws.Range("M132").Value = -2275.1

$ws.Range("H134").Value = 4523
$ws.Range("I134").Value = 1686.3
$ws.Range("J134").Value = 13978.667
$ws.Range("K134").Value = 5058.9
$ws.Range("L134").Value = 41936.001
$ws.Range("M134").Value = -2523.9
$ws.Range("N134").Value = -47006.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3823.5386
$ws.Range("J34").Value = 4121.3335
$ws.Range("L34").Value = 12364.0005
$ws.Range("N34").Value = -12532.0005

$ws.Range("H36").Value = 2157.1667
$ws.Range("I36").Value = 1988
$ws.Range("K36").Value = 5964
$ws.Range("M36").Value = -5795

$ws.Range("H39").Value = 8083.2104
$ws.Range("J39").Value = 8321.444
$ws.Range("L39").Value = 24964.332
$ws.Range("N39").Value = -25552.332

$ws.Range("H44").Value = 187.18182
$ws.Range("I44").Value = 187.18182
$ws.Range("K44").Value = 561.5454599999999
$ws.Range("M44").Value = -163.5454599999999

$ws.Range("H55").Value = 7260.393
$ws.Range("J55").Value = 7260.393
$ws.Range("L55").Value = 21781.179
$ws.Range("N55").Value = -22135.179

$ws.Range("H68").Value = 1042.1
$ws.Range("I68").Value = 707
$ws.Range("J68").Value = 1265.5
$ws.Range("K68").Value = 2121
$ws.Range("L68").Value = 3796.5
$ws.Range("M68").Value = -1310
$ws.Range("N68").Value = -5418.5

$ws.Range("H71").Value = 1042.1
$ws.Range("I71").Value = 707
$ws.Range("J71").Value = 1265.5
$ws.Range("K71").Value = 6363
$ws.Range("L71").Value = 11389.5
$ws.Range("M71").Value = -2307
$ws.Range("N71").Value = -19501.5

$ws.Range("H92").Value = 14250
$ws.Range("J92").Value = 14250
$ws.Range("L92").Value = 42750
$ws.Range("N92").Value = -45246

$ws.Range("H97").Value = 862.3333
$ws.Range("J97").Value = 1288
$ws.Range("L97").Value = 3864
$ws.Range("N97").Value = -4856

$ws.Range("H122").Value = 6798.1665
$ws.Range("I122").Value = 401.5
$ws.Range("J122").Value = 9996.5
$ws.Range("K122").Value = 3613.5
$ws.Range("L122").Value = 89968.5
$ws.Range("M122").Value = -1163.5
$ws.Range("N122").Value = -94868.5

$ws.Range("H132").Value = 1977.5555
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 59830.168
$ws.Range("J57").Value = 59830.168
$ws.Range("L57").Value = 59830.168
$ws.Range("N57").Value = -61470.168

$ws.Range("H80").Value = 3587.6667
$ws.Range("J80").Value = 7533
$ws.Range("L80").Value = 7533
$ws.Range("N80").Value = -9529

$ws.Range("H83").Value = 3587.6667
$ws.Range("J83").Value = 7533
$ws.Range("L83").Value = 37665
$ws.Range("N83").Value = -47649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5398.8823
$ws.Range("I46").Value = 4828.4
$ws.Range("K46").Value = 4828.4
$ws.Range("M46").Value = -4640.4

$ws.Range("H100").Value = 7177.778
$ws.Range("I100").Value = 4550
$ws.Range("J100").Value = 7928.5713
$ws.Range("K100").Value = 4550
$ws.Range("L100").Value = 7928.5713
$ws.Range("M100").Value = -4009
$ws.Range("N100").Value = -9010.5713

$ws.Range("H122").Value = 5845.4443
$ws.Range("I122").Value = 6012.8823
$ws.Range("K122").Value = 18038.6469
$ws.Range("M122").Value = -15588.6469

$ws.Range("H132").Value = 2158.0476
$ws.Range("I132").Value = 1850.7273
$ws.Range("K132").Value = 5552.1819
$ws.Range("M132").Value = -3022.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 17249.75
$ws.Range("I61").Value = 6333.3335
$ws.Range("J61").Value = 49999
$ws.Range("K61").Value = 6333.3335
$ws.Range("L61").Value = 49999
$ws.Range("M61").Value = -6041.3335
$ws.Range("N61").Value = -50583

$ws.Range("H81").Value = 3611.5715
$ws.Range("I81").Value = 3696.2
$ws.Range("K81").Value = 7392.4
$ws.Range("M81").Value = -6331.4

$ws.Range("H84").Value = 3611.5715
$ws.Range("I84").Value = 3696.2
$ws.Range("K84").Value = 36962
$ws.Range("M84").Value = -31658

$ws.Range("H122").Value = 3069.1785
$ws.Range("I122").Value = 2767.476
$ws.Range("J122").Value = 3974.2856
$ws.Range("K122").Value = 8302.428
$ws.Range("L122").Value = 11922.8568
$ws.Range("M122").Value = -5852.428
$ws.Range("N122").Value = -16822.8568

$ws.Range("H132").Value = 1663.014
$ws.Range("I132").Value = 1738.6833
$ws.Range("J132").Value = 1250.2727
$ws.Range("K132").Value = 5216.0499
$ws.Range("L132").Value = 3750.8181
$ws.Range("M132").Value = -2686.0499
$ws.Range("N132").Value = -8810.8181
